# Adds daily COVID bulletin rows 404-435 (2021-08-09 .. 2021-09-09) to Planilha1,
# matching the source commit "add data until September 9, 2021".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Columns A..J are raw daily counters; K..P are day-over-day deltas computed
# with formulas identical in shape to the ones already present in the sheet.
$data = @(
    @(404, 44417, 14946, 45,  6606, 21597, 6384, 77, 8, 69,  145),
    @(405, 44418, 14967, 68,  6616, 21651, 6387, 84, 8, 76,  145),
    @(406, 44419, 15041, 87,  6630, 21758, 6388, 97, 8, 89,  145),
    @(407, 44420, 15088, 80,  6635, 21803, 6406, 84, 6, 78,  145),
    @(408, 44421, 15130, 105, 6653, 21888, 6419, 89, 6, 83,  145),
    @(409, 44422, 15196, 47,  6675, 21918, 6430, 99, 5, 94,  146),
    @(410, 44423, 15196, 47,  6675, 21918, 6444, 85, 5, 80,  146),
    @(411, 44424, 15221, 62,  6680, 21963, 6457, 77, 5, 72,  146),
    @(412, 44425, 15254, 71,  6680, 22005, 6459, 74, 5, 69,  147),
    @(413, 44426, 15298, 79,  6698, 22075, 6459, 92, 5, 87,  147),
    @(414, 44427, 15331, 91,  6706, 22128, 6472, 87, 5, 82,  147),
    @(415, 44428, 15369, 79,  6714, 22162, 6489, 76, 5, 71,  149),
    @(416, 44429, 15422, 34,  6722, 22178, 6501, 72, 5, 67,  149),
    @(417, 44430, 15446, 30,  6724, 22200, 6513, 62, 4, 58,  149),
    @(418, 44431, 15469, 54,  6730, 22253, 6534, 47, 4, 43,  149),
    @(419, 44432, 15503, 51,  6738, 22292, 6537, 52, 4, 48,  149),
    @(420, 44433, 15514, 70,  6749, 22333, 6537, 63, 3, 60,  149),
    @(421, 44434, 15584, 88,  6753, 22425, 6537, 67, 3, 64,  149),
    @(422, 44435, 15648, 49,  6763, 22460, 6550, 63, 3, 60,  150),
    @(423, 44436, 15687, 45,  6767, 22499, 6555, 62, 3, 59,  150),
    @(424, 44437, 15704, 26,  6771, 22501, 6566, 55, 3, 52,  150),
    @(425, 44438, 15723, 43,  6773, 22539, 6566, 57, 3, 54,  150),
    @(426, 44439, 15750, 52,  6780, 22582, 6572, 58, 3, 55,  150),
    @(427, 44440, 15773, 52,  6789, 22614, 6573, 66, 3, 63,  150),
    @(428, 44441, 15811, 54,  6793, 22658, 6589, 54, 3, 51,  150),
    @(429, 44442, 15854, 38,  6799, 22691, 6600, 49, 5, 44,  150),
    @(430, 44443, 15868, 33,  6805, 22706, 6608, 47, 5, 42,  150),
    @(431, 44444, 15893, 12,  6805, 22710, 6614, 41, 5, 36,  150),
    @(432, 44445, 15893, 22,  6806, 22721, 6620, 36, 5, 31,  150),
    @(433, 44446, 15894, 42,  6813, 22752, 6626, 37, 5, 32,  150),
    @(434, 44447, 15918, 44,  6816, 22778, 6626, 40, 5, 35,  150),
    @(435, 44448, 15928, 52,  6818, 22798, 6630, 37, 5, 32,  151)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value  = $row[1]   # A DATA
    $ws.Cells.Item($r, 2).Value  = $row[2]   # B CONFIRMADOS
    $ws.Cells.Item($r, 3).Value  = $row[3]   # C EM INVESTIGACAO
    $ws.Cells.Item($r, 4).Value  = $row[4]   # D OBITOS
    $ws.Cells.Item($r, 5).Value  = $row[5]   # E EXAMINADOS
    $ws.Cells.Item($r, 6).Value  = $row[6]   # F RECUPERADOS
    $ws.Cells.Item($r, 7).Value  = $row[7]   # G ATIVOS
    $ws.Cells.Item($r, 8).Value  = $row[8]   # H HOSPITAL
    $ws.Cells.Item($r, 9).Value  = $row[9]   # I DOMICILIO
    $ws.Cells.Item($r, 10).Value = $row[10]  # J VERSAO
}

# Day-over-day delta formulas, filled in the same two batches (404:422 then
# 423:435) and column order (K,L,M,N,O,P) that the original workbook used, so
# shared-formula grouping matches how the sheet was actually built up.
$ws.Range("K404:K422").Formula = "=D404-D403"
$ws.Range("L404:L422").Formula = "=F404-F403"
$ws.Range("M404:M422").Formula = "=B404-B403"
$ws.Range("N404:N422").Formula = "=J404-J403"
$ws.Range("O404:O422").Formula = "=G404-G403"
$ws.Range("P404:P422").Formula = "=C404-C403"

$ws.Range("K423:K435").Formula = "=D423-D422"
$ws.Range("L423:L435").Formula = "=F423-F422"
$ws.Range("M423:M435").Formula = "=B423-B422"
$ws.Range("N423:N435").Formula = "=J423-J422"
$ws.Range("O423:O435").Formula = "=G423-G422"
$ws.Range("P423:P435").Formula = "=C423-C422"

# Restore the header-row freeze and move the live selection to where the
# author last left off (mirrors the view-state change in the diff).
$win = $excel.ActiveWindow
$win.FreezePanes = $false
$ws.Range("A2").Select()
$win.FreezePanes = $true
$ws.Range("R416").Select()

Write-Output "Added rows 404:435"
